$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# --- Update the export-time value (B2): "2020-12-22" -> "2021-01-11" ---
# B2 keeps its original "Text"-like alignment style but is a plain text value,
# so force a Text number format while assigning, then restore the style's
# original number format so the cell keeps looking the same as before.
$b2 = $ws.Range("B2")
$b2Format = $b2.NumberFormat
$b2.NumberFormat = "@"
$b2.Value = "2021-01-11"
$b2.NumberFormat = $b2Format

# --- Replace the monthly data block D4:I15 with the new dataset ---
# These columns are formatted with a "Text" number format (numFmtId 49) even
# though they store real numbers, so we temporarily switch to General while
# writing the values to avoid Excel auto-converting them into text, then
# restore the original number format afterwards.
$dataRange = $ws.Range("D4:I15")
$origFormat = $dataRange.NumberFormat
$dataRange.NumberFormat = "General"

$newValues = @{
    4  = @(51.770907643488,    62.0, 3.779187243024677,  585.0, 29.918763317971507, 249.0)
    5  = @(16.1890817198666,   41.0, 57.2267239689793,   443.0, 42.81109579543078,    8.0)
    6  = @(34.865603386588404, 30.0, 5.930332835531042,  467.0, 19.14982628027355,  241.0)
    7  = @(62.96257363175009,  24.0, 143.50191413190825, 475.0, 42.258868348470855,  95.0)
    8  = @(19.889574198994808, 48.0, 149.75682073444975, 554.0, 22.46973573891287,  262.0)
    9  = @(8.652804050308568,   2.0, 119.30554915254275, 601.0, 7.364525953548196,  220.0)
    10 = @(2.093522821491483,   8.0, 54.29677905940273,  602.0, 29.657148134498776, 218.0)
    11 = @(47.97203229235001,  31.0, 94.25180927462719,  326.0, 17.71348947458766,  113.0)
    12 = @(7.6895749100910615, 38.0, 7.368976908940475,  525.0, 14.264113169648962, 177.0)
    13 = @(19.413483175636713, 31.0, 91.54756804917443,  242.0, 24.89281261926465,  222.0)
    14 = @(24.04298275497596,  41.0, 68.42584415275893,  435.0, 5.277532011778943,  229.0)
    15 = @(45.38692115455145,   8.0, 2.314368768090385,  317.0, 13.366202671283098,  26.0)
}

$columns = @("D", "E", "F", "G", "H", "I")
foreach ($row in $newValues.Keys) {
    $values = $newValues[$row]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $cellRef = "$($columns[$i])$row"
        $ws.Range($cellRef).Value = $values[$i]
    }
}

$dataRange.NumberFormat = $origFormat
